# Auto-generated edit script applying numeric corrections to Sheets per commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 146
$ws.Range("I12").Value = 146
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 146
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 24
$ws.Range("N12").ClearContents()
$ws.Range("H43").Value = 200
$ws.Range("J43").Value = 200
$ws.Range("L43").Value = 200
$ws.Range("N43").Value = -338
$ws.Range("H51").Value = 3060.08
$ws.Range("J51").Value = 2961.6924
$ws.Range("L51").Value = 2961.6924
$ws.Range("N51").Value = -3929.6924
$ws.Range("H86").Value = 62505116
$ws.Range("I86").Value = 5659.143
$ws.Range("J86").Value = 111115800
$ws.Range("K86").Value = 5659.143
$ws.Range("L86").Value = 111115800
$ws.Range("M86").Value = -4536.143
$ws.Range("N86").Value = -111118046
$ws.Range("H87").Value = 173435
$ws.Range("J87").Value = 173435
$ws.Range("L87").Value = 173435
$ws.Range("N87").Value = -175931
$ws.Range("H89").Value = 62505116
$ws.Range("I89").Value = 5659.143
$ws.Range("J89").Value = 111115800
$ws.Range("K89").Value = 28295.715
$ws.Range("L89").Value = 555579000
$ws.Range("M89").Value = -22679.715
$ws.Range("N89").Value = -555590232
$ws.Range("H90").Value = 173435
$ws.Range("J90").Value = 173435
$ws.Range("L90").Value = 520305
$ws.Range("N90").Value = -532785
$ws.Range("H99").Value = 20833964
$ws.Range("J99").Value = 1194.5
$ws.Range("L99").Value = 3583.5
$ws.Range("N99").Value = -6579.5
$ws.Range("H108").Value = 98888.75
$ws.Range("J108").Value = 98888.75
$ws.Range("L108").Value = 98888.75
$ws.Range("N108").Value = -106568.75
$ws.Range("H114").Value = 58717.5
$ws.Range("J114").Value = 58717.5
$ws.Range("L114").Value = 58717.5
$ws.Range("N114").Value = -67395.5
$ws.Range("H137").Value = 3467.3704
$ws.Range("I137").Value = 2592.1738
$ws.Range("K137").Value = 7776.5214
$ws.Range("M137").Value = -5226.5214
$ws.Range("H138").Value = 2890.4082
$ws.Range("I138").Value = 2018.3182
$ws.Range("K138").Value = 6054.9546
$ws.Range("M138").Value = -914.9546
$ws.Range("H140").Value = 207855.28
$ws.Range("I140").Value = 64995
$ws.Range("J140").Value = 231665.33
$ws.Range("K140").Value = 64995
$ws.Range("L140").Value = 231665.33
$ws.Range("M140").Value = -59815
$ws.Range("N140").Value = -242025.33
$ws.Range("H141").Value = 5900.8
$ws.Range("I141").Value = 4889.778
$ws.Range("J141").Value = 15000
$ws.Range("K141").Value = 14669.334
$ws.Range("L141").Value = 45000
$ws.Range("M141").Value = -9489.334000000001
$ws.Range("N141").Value = -55360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1416.75
$ws.Range("I2").Value = 804.04
$ws.Range("K2").Value = 804.04
$ws.Range("M2").Value = -691.04
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H32").Value = 18944562
$ws.Range("I32").Value = 15154761
$ws.Range("J32").Value = 22734364
$ws.Range("K32").Value = 15154761
$ws.Range("L32").Value = 22734364
$ws.Range("M32").Value = -15154474
$ws.Range("N32").Value = -22734938
$ws.Range("H52").Value = 112494.5
$ws.Range("J52").Value = 112494.5
$ws.Range("L52").Value = 112494.5
$ws.Range("N52").Value = -113130.5
$ws.Range("H61").Value = 2635.3635
$ws.Range("I61").Value = 2648.9
$ws.Range("K61").Value = 2648.9
$ws.Range("M61").Value = -2436.9
$ws.Range("H97").Value = 812.3226
$ws.Range("I97").Value = 614.8889
$ws.Range("K97").Value = 614.8889
$ws.Range("M97").Value = -118.8889
$ws.Range("H102").Value = 2538
$ws.Range("I102").Value = 2208
$ws.Range("J102").Value = 3198
$ws.Range("K102").Value = 2208
$ws.Range("L102").Value = 3198
$ws.Range("M102").Value = -586
$ws.Range("N102").Value = -6442
$ws.Range("H116").Value = 1416.75
$ws.Range("I116").Value = 804.04
$ws.Range("K116").Value = 804.04
$ws.Range("M116").Value = 1489.96
$ws.Range("H122").Value = 2098.7568
$ws.Range("I122").Value = 1833
$ws.Range("K122").Value = 5499
$ws.Range("M122").Value = -3049
$ws.Range("H132").Value = 2969.158
$ws.Range("I132").Value = 2720.08
$ws.Range("J132").Value = 3448.1538
$ws.Range("K132").Value = 8160.24
$ws.Range("L132").Value = 10344.4614
$ws.Range("M132").Value = -5630.24
$ws.Range("N132").Value = -15404.4614
$ws.Range("H136").Value = 2635.3635
$ws.Range("I136").Value = 2648.9
$ws.Range("K136").Value = 7946.700000000001
$ws.Range("M136").Value = -5396.700000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1416.75
$ws.Range("I3").Value = 804.04
$ws.Range("K3").Value = 804.04
$ws.Range("M3").Value = -690.04
$ws.Range("H29").Value = 166.66667
$ws.Range("I29").Value = 166.66667
$ws.Range("K29").Value = 166.66667
$ws.Range("M29").Value = 122.33333
$ws.Range("H105").Value = 2834.9656
$ws.Range("I105").Value = 2245.842
$ws.Range("K105").Value = 2245.842
$ws.Range("M105").Value = -498.8420000000001
$ws.Range("H134").Value = 10941553
$ws.Range("I134").Value = 2382404
$ws.Range("K134").Value = 7147212
$ws.Range("M134").Value = -7144677

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1028.7
$ws.Range("I16").Value = 930.1667
$ws.Range("J16").Value = 1176.5
$ws.Range("K16").Value = 930.1667
$ws.Range("L16").Value = 1176.5
$ws.Range("M16").Value = -643.1667
$ws.Range("N16").Value = -1750.5
$ws.Range("H18").Value = 40591
$ws.Range("J18").Value = 40591
$ws.Range("L18").Value = 40591
$ws.Range("N18").Value = -41051
$ws.Range("H31").Value = 1405356.9
$ws.Range("I31").Value = 2809.6667
$ws.Range("K31").Value = 2809.6667
$ws.Range("M31").Value = -2514.6667
$ws.Range("H34").Value = 1405356.9
$ws.Range("I34").Value = 2809.6667
$ws.Range("K34").Value = 2809.6667
$ws.Range("M34").Value = -2607.6667
$ws.Range("H50").Value = 39250
$ws.Range("H53").Value = 36089.8
$ws.Range("J53").Value = 36089.8
$ws.Range("L53").Value = 36089.8
$ws.Range("N53").Value = -37303.8
$ws.Range("H58").Value = 2705.0588
$ws.Range("I58").Value = 2164
$ws.Range("J58").Value = 3697
$ws.Range("K58").Value = 2164
$ws.Range("L58").Value = 3697
$ws.Range("M58").Value = -1961
$ws.Range("N58").Value = -4103
$ws.Range("H107").Value = 51856.3
$ws.Range("I107").Value = 77877.234
$ws.Range("K107").Value = 77877.234
$ws.Range("M107").Value = -75957.234
$ws.Range("H108").Value = 29879
$ws.Range("J108").Value = 29879
$ws.Range("L108").Value = 29879
$ws.Range("N108").Value = -37559
$ws.Range("H111").Value = 48970
$ws.Range("J111").Value = 48970
$ws.Range("L111").Value = 48970
$ws.Range("N111").Value = -57150
$ws.Range("H112").Value = 145663
$ws.Range("J112").Value = 145663
$ws.Range("L112").Value = 145663
$ws.Range("N112").Value = -148617
$ws.Range("H113").Value = 1028.7
$ws.Range("I113").Value = 930.1667
$ws.Range("J113").Value = 1176.5
$ws.Range("K113").Value = 930.1667
$ws.Range("L113").Value = 1176.5
$ws.Range("M113").Value = 1239.8333
$ws.Range("N113").Value = -5516.5
$ws.Range("H114").Value = 88791.5
$ws.Range("J114").Value = 88791.5
$ws.Range("L114").Value = 88791.5
$ws.Range("N114").Value = -97469.5
$ws.Range("H116").Value = 105993
$ws.Range("J116").Value = 105993
$ws.Range("L116").Value = 105993
$ws.Range("N116").Value = -115171
$ws.Range("H118").Value = 139990
$ws.Range("J118").Value = 139990
$ws.Range("L118").Value = 139990
$ws.Range("N118").Value = -143304
$ws.Range("H122").Value = 1822.2174
$ws.Range("J122").Value = 2299.8572
$ws.Range("L122").Value = 6899.571599999999
$ws.Range("N122").Value = -11799.5716
$ws.Range("H132").Value = 2086.6365
$ws.Range("I132").Value = 899.2857
$ws.Range("K132").Value = 2697.8571
$ws.Range("M132").Value = -167.8571000000002
$ws.Range("H136").Value = 2705.0588
$ws.Range("I136").Value = 2164
$ws.Range("J136").Value = 3697
$ws.Range("K136").Value = 6492
$ws.Range("L136").Value = 11091
$ws.Range("M136").Value = -3942
$ws.Range("N136").Value = -16191
$ws.Range("H141").Value = 938131.5
$ws.Range("J141").Value = 1154748.8
$ws.Range("L141").Value = 1154748.8
$ws.Range("N141").Value = -1165108.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 34370268
$ws.Range("I4").Value = 30010842
$ws.Range("K4").Value = 90032526
$ws.Range("M4").Value = -90032414
$ws.Range("H118").Value = 1671.1666
$ws.Range("I118").Value = 1205.4
$ws.Range("K118").Value = 3616.2
$ws.Range("M118").Value = -2373.2
$ws.Range("H131").Value = 1657.7142
$ws.Range("I131").Value = 1318.2222
$ws.Range("K131").Value = 3954.6666
$ws.Range("M131").Value = 1085.3334
$ws.Range("H132").Value = 1495.875
$ws.Range("I132").Value = 539.5
$ws.Range("K132").Value = 4855.5
$ws.Range("M132").Value = -2325.5
$ws.Range("H137").Value = 9686.666999999999
$ws.Range("I137").Value = 1508.8334
$ws.Range("J137").Value = 17864.5
$ws.Range("K137").Value = 4526.5002
$ws.Range("L137").Value = 53593.5
$ws.Range("M137").Value = 573.4997999999996
$ws.Range("N137").Value = -63793.5
$ws.Range("H141").Value = 6333.222
$ws.Range("I141").Value = 1999.8334
$ws.Range("J141").Value = 15000
$ws.Range("K141").Value = 5999.5002
$ws.Range("L141").Value = 45000
$ws.Range("M141").Value = -819.5002000000004
$ws.Range("N141").Value = -55360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 25071.857
$ws.Range("I20").Value = 11751.75
$ws.Range("J20").Value = 42832
$ws.Range("K20").Value = 11751.75
$ws.Range("L20").Value = 42832
$ws.Range("M20").Value = -11506.75
$ws.Range("N20").Value = -43322
$ws.Range("H24").Value = 41000.75
$ws.Range("I24").Value = 24006
$ws.Range("J24").Value = 46665.668
$ws.Range("K24").Value = 24006
$ws.Range("L24").Value = 46665.668
$ws.Range("M24").Value = -23833
$ws.Range("N24").Value = -47011.668
$ws.Range("H70").Value = 4443.8887
$ws.Range("I70").Value = 3995
$ws.Range("K70").Value = 3995
$ws.Range("M70").Value = -3725
$ws.Range("H73").Value = 4443.8887
$ws.Range("I73").Value = 3995
$ws.Range("K73").Value = 3995
$ws.Range("M73").Value = -3059
$ws.Range("H80").Value = 2343.375
$ws.Range("J80").Value = 3499.6667
$ws.Range("L80").Value = 3499.6667
$ws.Range("N80").Value = -5495.6667
$ws.Range("H83").Value = 2343.375
$ws.Range("J83").Value = 3499.6667
$ws.Range("L83").Value = 17498.3335
$ws.Range("N83").Value = -27482.3335
$ws.Range("H97").Value = 831.03705
$ws.Range("I97").Value = 469
$ws.Range("K97").Value = 469
$ws.Range("M97").Value = 27
$ws.Range("H122").Value = 1722.75
$ws.Range("J122").Value = 1344.2222
$ws.Range("L122").Value = 4032.6666
$ws.Range("N122").Value = -8932.6666
$ws.Range("H126").Value = 4659.963
$ws.Range("I126").Value = 4255
$ws.Range("K126").Value = 12765
$ws.Range("M126").Value = -10295

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 17070.08
$ws.Range("I7").Value = 16837.6
$ws.Range("J7").Value = 18000
$ws.Range("K7").Value = 16837.6
$ws.Range("L7").Value = 18000
$ws.Range("M7").Value = -16725.6
$ws.Range("N7").Value = -18224
$ws.Range("H16").Value = 1760.8667
$ws.Range("J16").Value = 3000
$ws.Range("L16").Value = 3000
$ws.Range("N16").Value = -3340
$ws.Range("H40").Value = 8424.368
$ws.Range("I40").Value = 8239.058999999999
$ws.Range("K40").Value = 8239.058999999999
$ws.Range("M40").Value = -8103.058999999999
$ws.Range("H55").Value = 245.05882
$ws.Range("I55").Value = 167.35
$ws.Range("J55").Value = 356.07144
$ws.Range("K55").Value = 167.35
$ws.Range("L55").Value = 356.07144
$ws.Range("M55").Value = 5.650000000000006
$ws.Range("N55").Value = -702.0714399999999
$ws.Range("H93").Value = 1085.5714
$ws.Range("I93").Value = 837.2727
$ws.Range("K93").Value = 837.2727
$ws.Range("M93").Value = 410.7273
$ws.Range("H126").Value = 17070.08
$ws.Range("I126").Value = 16837.6
$ws.Range("J126").Value = 18000
$ws.Range("K126").Value = 50512.8
$ws.Range("L126").Value = 54000
$ws.Range("M126").Value = -48042.8
$ws.Range("N126").Value = -58940
$ws.Range("H132").Value = 5142.5713
$ws.Range("I132").Value = 4666.3335
$ws.Range("J132").Value = 5499.75
$ws.Range("K132").Value = 13999.0005
$ws.Range("L132").Value = 16499.25
$ws.Range("M132").Value = -11469.0005
$ws.Range("N132").Value = -21559.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 7002200
$ws.Range("I2").Value = 7694827
$ws.Range("J2").Value = 2500125
$ws.Range("K2").Value = 7694827
$ws.Range("L2").Value = 2500125
$ws.Range("M2").Value = -7694715
$ws.Range("N2").Value = -2500349
$ws.Range("H21").Value = 25012.25
$ws.Range("I21").Value = 15015
$ws.Range("J21").Value = 28344.666
$ws.Range("K21").Value = 15015
$ws.Range("L21").Value = 28344.666
$ws.Range("M21").Value = -14780
$ws.Range("N21").Value = -28814.666
$ws.Range("H22").Value = 11405.8
$ws.Range("I22").Value = 5000
$ws.Range("J22").Value = 13007.25
$ws.Range("K22").Value = 5000
$ws.Range("L22").Value = 13007.25
$ws.Range("M22").Value = -4707
$ws.Range("N22").Value = -13593.25
$ws.Range("H35").Value = 25012.25
$ws.Range("I35").Value = 15015
$ws.Range("J35").Value = 28344.666
$ws.Range("K35").Value = 15015
$ws.Range("L35").Value = 28344.666
$ws.Range("M35").Value = -14725
$ws.Range("N35").Value = -28924.666
$ws.Range("H100").Value = 2415.875
$ws.Range("I100").Value = 2171.1667
$ws.Range("J100").Value = 3150
$ws.Range("K100").Value = 4342.3334
$ws.Range("L100").Value = 6300
$ws.Range("M100").Value = -3801.3334
$ws.Range("N100").Value = -7382
$ws.Range("H115").Value = 60377
$ws.Range("J115").Value = 60377
$ws.Range("L115").Value = 60377
$ws.Range("N115").Value = -63511
$ws.Range("H116").Value = 59995.5
$ws.Range("J116").Value = 59995.5
$ws.Range("L116").Value = 59995.5
$ws.Range("N116").Value = -69173.5
$ws.Range("H118").Value = 107999
$ws.Range("J118").Value = 107999
$ws.Range("L118").Value = 107999
$ws.Range("N118").Value = -111313
$ws.Range("H119").Value = 27500
$ws.Range("J119").Value = 25000
$ws.Range("L119").Value = 25000
$ws.Range("N119").Value = -34676
$ws.Range("H121").Value = 62644
$ws.Range("J121").Value = 62644
$ws.Range("L121").Value = 62644
$ws.Range("N121").Value = -66138
$ws.Range("H132").Value = 2513.8057
$ws.Range("I132").Value = 1818.1724
$ws.Range("K132").Value = 5454.5172
$ws.Range("M132").Value = -2924.5172
$ws.Range("H135").Value = 99299.664
$ws.Range("J135").Value = 99299.664
$ws.Range("L135").Value = 99299.664
$ws.Range("N135").Value = -109439.664
$ws.Range("H136").Value = 22420.607
$ws.Range("I136").Value = 1726.8214
$ws.Range("K136").Value = 5180.4642
$ws.Range("M136").Value = -2630.4642
